# btmigrate_work.xlsx - smart rule additions: merge in new rule rows, add E2 "winscp"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a numeric-looking string (e.g. "1453") as real text, not a number,
# by building it through a text formula in a scratch cell and pasting values only.
# This avoids Excel's normal behaviour of coercing digit-only strings to numbers.
function Set-TextValue($addr, $val) {
    $ws.Range("ZZ1").Formula = '="' + $val + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# Row 8: new "admin" smart rule (oracle) -- set before E2 so new shared
# strings are appended in the same order as the target workbook ("admin" first)
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "admin"
$ws.Range("C8").Value = "10.59.76.111"
$ws.Range("D8").Value = "tkracdb.thy.com"
$ws.Range("E8").Value = "oracle"
$ws.Range("F8").Value = "oracle"
$ws.Range("G8").Value = "E_MISIR"
$ws.Range("I8").Value = "tkrac"
Set-TextValue "J8" "1453"
$ws.Range("K8").Value = "oracle"

# New application value for the already-existing first rule row (E2)
$ws.Range("E2").Value = "winscp"

# Row 9
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "pam118064"
$ws.Range("C9").Value = "10.59.76.112"
$ws.Range("D9").Value = "deneme.thynet.thy.com"
$ws.Range("E9").Value = "mssql"
$ws.Range("F9").Value = "mssql"
$ws.Range("G9").Value = "S_OZCAN"
$ws.Range("H9").Value = "S_OZCAN"
$ws.Range("I9").Value = "nan"
Set-TextValue "J9" "1433"
$ws.Range("K9").Value = "mssql"
$ws.Range("L9").Value = "quasys.local"

# Row 10
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "pam100887"
$ws.Range("C10").Value = "10.58.59.101"
$ws.Range("D10").Value = "host-002"
$ws.Range("F10").Value = "Windows"
$ws.Range("G10").Value = "E_PEKDAS"
$ws.Range("I10").Value = "nan"
$ws.Range("J10").Value = "nan"
$ws.Range("K10").Value = "domain"
$ws.Range("L10").Value = "quasys.local"

# Row 11
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "pam100888"
$ws.Range("C11").Value = "10.58.59.102"
$ws.Range("D11").Value = "host-010"
$ws.Range("F11").Value = "Windows"
$ws.Range("G11").Value = "E_YARDIMCI"
$ws.Range("H11").Value = "E_YARDIMCI"
$ws.Range("I11").Value = "nan"
$ws.Range("J11").Value = "nan"
$ws.Range("K11").Value = "domain"
$ws.Range("L11").Value = "quasys.local"

# Row 12
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "PAM100545"
$ws.Range("C12").Value = "10.59.76.111"
$ws.Range("D12").Value = "tkracdb.thy.com"
$ws.Range("E12").Value = "oracle"
$ws.Range("F12").Value = "oracle"
$ws.Range("G12").Value = "S_OZCAN"
$ws.Range("I12").Value = "tkrac"
Set-TextValue "J12" "1453"
$ws.Range("K12").Value = "oracle"

# Row 13
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "pam100888"
$ws.Range("C13").Value = "10.59.76.111"
$ws.Range("D13").Value = "tkracdb.thy.com"
$ws.Range("E13").Value = "oracle"
$ws.Range("F13").Value = "oracle"
$ws.Range("G13").Value = "E_PEKDAS"
$ws.Range("I13").Value = "tkrac"
Set-TextValue "J13" "1453"
$ws.Range("K13").Value = "oracle"

# Clean up the scratch cell used for text coercion
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = 0

# Match the selection state left behind in the saved file
$ws.Range("I17").Select()
